$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 363
$ws.Range("I4").Value = 188.15384
$ws.Range("J4").Value = 1499.5
$ws.Range("K4").Value = 188.15384
$ws.Range("L4").Value = 1499.5
$ws.Range("M4").Value = -74.15384
$ws.Range("N4").Value = -1727.5
$ws.Range("H6").Value = 1014.4286
$ws.Range("J6").Value = 1300.6666
$ws.Range("L6").Value = 3901.9998
$ws.Range("N6").Value = -4125.9998
$ws.Range("H33").Value = 374.72223
$ws.Range("I33").Value = 396.91666
$ws.Range("J33").Value = 330.33334
$ws.Range("K33").Value = 396.91666
$ws.Range("L33").Value = 330.33334
$ws.Range("M33").Value = -167.91666
$ws.Range("N33").Value = -788.33334
$ws.Range("H41").Value = 333.5
$ws.Range("I41").Value = 379.2353
$ws.Range("J41").Value = 222.42857
$ws.Range("K41").Value = 379.2353
$ws.Range("L41").Value = 222.42857
$ws.Range("M41").Value = 60.7647
$ws.Range("N41").Value = -1102.42857
$ws.Range("H43").Value = 15500.2
$ws.Range("J43").Value = 13000.333
$ws.Range("L43").Value = 13000.333
$ws.Range("N43").Value = -13138.333
$ws.Range("H99").Value = 793
$ws.Range("I99").Value = 589.5
$ws.Range("J99").Value = 1200
$ws.Range("K99").Value = 1768.5
$ws.Range("L99").Value = 3600
$ws.Range("M99").Value = -270.5
$ws.Range("N99").Value = -6596
$ws.Range("H103").Value = 580.5833
$ws.Range("I103").Value = 466.7143
$ws.Range("J103").Value = 740
$ws.Range("K103").Value = 1400.1429
$ws.Range("L103").Value = 2220
$ws.Range("M103").Value = -814.1428999999998
$ws.Range("N103").Value = -3392
$ws.Range("H125").Value = 7468.55
$ws.Range("I125").Value = 7252.5454
$ws.Range("J125").Value = 7732.5557
$ws.Range("K125").Value = 65272.9086
$ws.Range("L125").Value = 69593.0013
$ws.Range("M125").Value = -62812.9086
$ws.Range("N125").Value = -74513.0013
$ws.Range("H138").Value = 5115.6113
$ws.Range("J138").Value = 5595.2896
$ws.Range("L138").Value = 16785.8688
$ws.Range("N138").Value = -27065.8688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 7
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 105
$ws.Range("N5").Value = ""
$ws.Range("H32").Value = 24595.572
$ws.Range("I32").Value = 24595.572
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 24595.572
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -24308.572
$ws.Range("N32").Value = ""
$ws.Range("H45").Value = 3200
$ws.Range("I45").Value = 3500
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 3500
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -3123
$ws.Range("N45").Value = -2754
$ws.Range("H110").Value = 9261399
$ws.Range("J110").Value = 3901.2
$ws.Range("L110").Value = 3901.2
$ws.Range("N110").Value = -7991.2
$ws.Range("H112").Value = 99997.5
$ws.Range("J112").Value = 99997.5
$ws.Range("L112").Value = 99997.5
$ws.Range("N112").Value = -102951.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 7
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 7
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 108
$ws.Range("N4").Value = ""
$ws.Range("H134").Value = 4502.344
$ws.Range("I134").Value = 3293.2546
$ws.Range("J134").Value = 15585.667
$ws.Range("K134").Value = 9879.763800000001
$ws.Range("L134").Value = 46757.001
$ws.Range("M134").Value = -7344.763800000001
$ws.Range("N134").Value = -51827.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 814.4286
$ws.Range("J2").Value = 1350
$ws.Range("L2").Value = 1350
$ws.Range("N2").Value = -1576
$ws.Range("H31").Value = 34490444
$ws.Range("I31").Value = 142861920
$ws.Range("J31").Value = 8607.409
$ws.Range("K31").Value = 142861920
$ws.Range("L31").Value = 8607.409
$ws.Range("M31").Value = -142861625
$ws.Range("N31").Value = -9197.409
$ws.Range("H34").Value = 34490444
$ws.Range("I34").Value = 142861920
$ws.Range("J34").Value = 8607.409
$ws.Range("K34").Value = 142861920
$ws.Range("L34").Value = 8607.409
$ws.Range("M34").Value = -142861718
$ws.Range("N34").Value = -9011.409
$ws.Range("H105").Value = 1923.8334
$ws.Range("I105").Value = 1565.1111
$ws.Range("K105").Value = 1565.1111
$ws.Range("M105").Value = 181.8888999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 27818446
$ws.Range("I4").Value = 37884400
$ws.Range("J4").Value = 5363632
$ws.Range("K4").Value = 113653200
$ws.Range("L4").Value = 16090896
$ws.Range("M4").Value = -113653088
$ws.Range("N4").Value = -16091120
$ws.Range("H136").Value = 5424
$ws.Range("I136").Value = 2815
$ws.Range("J136").Value = 8033
$ws.Range("K136").Value = 8445
$ws.Range("L136").Value = 24099
$ws.Range("M136").Value = -3345
$ws.Range("N136").Value = -34299

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7868.778
$ws.Range("I70").Value = 6469
$ws.Range("K70").Value = 6469
$ws.Range("M70").Value = -6199
$ws.Range("H73").Value = 7868.778
$ws.Range("I73").Value = 6469
$ws.Range("K73").Value = 6469
$ws.Range("M73").Value = -5533
$ws.Range("H80").Value = 4386.4614
$ws.Range("I80").Value = 3845.4
$ws.Range("J80").Value = 4724.625
$ws.Range("K80").Value = 3845.4
$ws.Range("L80").Value = 4724.625
$ws.Range("M80").Value = -2847.4
$ws.Range("N80").Value = -6720.625
$ws.Range("H83").Value = 4386.4614
$ws.Range("I83").Value = 3845.4
$ws.Range("J83").Value = 4724.625
$ws.Range("K83").Value = 19227
$ws.Range("L83").Value = 23623.125
$ws.Range("M83").Value = -14235
$ws.Range("N83").Value = -33607.125
$ws.Range("H99").Value = 9752.909
$ws.Range("I99").Value = 9752.909
$ws.Range("M99").Value = -7506.909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1175.7941
$ws.Range("I16").Value = 914.48
$ws.Range("J16").Value = 1901.6666
$ws.Range("K16").Value = 914.48
$ws.Range("L16").Value = 1901.6666
$ws.Range("M16").Value = -744.48
$ws.Range("N16").Value = -2241.6666
$ws.Range("H46").Value = 7147.2
$ws.Range("I46").Value = 1110.25
$ws.Range("K46").Value = 1110.25
$ws.Range("M46").Value = -922.25
$ws.Range("H55").Value = 544.3333
$ws.Range("I55").Value = 305.47058
$ws.Range("K55").Value = 305.47058
$ws.Range("M55").Value = -132.47058
$ws.Range("H68").Value = 5462.5
$ws.Range("I68").Value = 2666.6667
$ws.Range("K68").Value = 2666.6667
$ws.Range("M68").Value = -1917.6667
$ws.Range("H71").Value = 5462.5
$ws.Range("I71").Value = 2666.6667
$ws.Range("K71").Value = 13333.3335
$ws.Range("M71").Value = -9589.333500000001
$ws.Range("H110").Value = 39997.5
$ws.Range("J110").Value = 39997.5
$ws.Range("L110").Value = 39997.5
$ws.Range("N110").Value = -48177.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 6300.143
$ws.Range("J41").Value = 6300.143
$ws.Range("L41").Value = 6300.143
$ws.Range("N41").Value = -7080.143
$ws.Range("H81").Value = 12890.243
$ws.Range("I81").Value = 4078.4167
$ws.Range("K81").Value = 8156.8334
$ws.Range("M81").Value = -7095.8334
$ws.Range("H84").Value = 12890.243
$ws.Range("I84").Value = 4078.4167
$ws.Range("K84").Value = 40784.167
$ws.Range("M84").Value = -35480.167
$ws.Range("H96").Value = 4585.3
$ws.Range("I96").Value = 3231.625
$ws.Range("K96").Value = 3231.625
$ws.Range("M96").Value = -1858.625
$ws.Range("H100").Value = 2094.476
$ws.Range("I100").Value = 1198.7273
$ws.Range("J100").Value = 3079.8
$ws.Range("K100").Value = 2397.4546
$ws.Range("L100").Value = 6159.6
$ws.Range("M100").Value = -1856.4546
$ws.Range("N100").Value = -7241.6
$ws.Range("H106").Value = 75000
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").Value = ""
